$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column ("D") values are plain numeric-looking text (e.g. "386.00",
# "0.0770") in the source data. Assigning such a string straight to .Value
# lets Excel auto-detect it as a number, silently dropping the trailing
# zero / precision (e.g. "386.00" -> 386, "6.90" -> 6.9). To keep the exact
# literal text, format the cell as Text first, then restore the "Normal"
# style afterwards so no stray number-format style lingers on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.354.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.375.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.19%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.375.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.950.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.377.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.468.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -7.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.518.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.07%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.405.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0770"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.778"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.519.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  -3.13%  "
